$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap F:V data between paired rows (home/away ordering corrections) ---
$rowA = $ws.Range("F63:V63").Value2
$rowB = $ws.Range("F64:V64").Value2
$ws.Range("F63:V63").Value2 = $rowB
$ws.Range("F64:V64").Value2 = $rowA

$rowA = $ws.Range("F68:V68").Value2
$rowB = $ws.Range("F69:V69").Value2
$ws.Range("F68:V68").Value2 = $rowB
$ws.Range("F69:V69").Value2 = $rowA

$rowA = $ws.Range("F79:V79").Value2
$rowB = $ws.Range("F80:V80").Value2
$ws.Range("F79:V79").Value2 = $rowB
$ws.Range("F80:V80").Value2 = $rowA

$rowA = $ws.Range("F87:V87").Value2
$rowB = $ws.Range("F88:V88").Value2
$ws.Range("F87:V87").Value2 = $rowB
$ws.Range("F88:V88").Value2 = $rowA

$rowA = $ws.Range("F89:V89").Value2
$rowB = $ws.Range("F90:V90").Value2
$ws.Range("F89:V89").Value2 = $rowB
$ws.Range("F90:V90").Value2 = $rowA

$rowA = $ws.Range("F91:V91").Value2
$rowB = $ws.Range("F92:V92").Value2
$ws.Range("F91:V91").Value2 = $rowB
$ws.Range("F92:V92").Value2 = $rowA

$rowA = $ws.Range("F94:V94").Value2
$rowB = $ws.Range("F95:V95").Value2
$ws.Range("F94:V94").Value2 = $rowB
$ws.Range("F95:V95").Value2 = $rowA

$rowA = $ws.Range("F98:V98").Value2
$rowB = $ws.Range("F99:V99").Value2
$ws.Range("F98:V98").Value2 = $rowB
$ws.Range("F99:V99").Value2 = $rowA

$rowA = $ws.Range("F106:V106").Value2
$rowB = $ws.Range("F107:V107").Value2
$ws.Range("F106:V106").Value2 = $rowB
$ws.Range("F107:V107").Value2 = $rowA

$rowA = $ws.Range("F112:V112").Value2
$rowB = $ws.Range("F113:V113").Value2
$ws.Range("F112:V112").Value2 = $rowB
$ws.Range("F113:V113").Value2 = $rowA

$rowA = $ws.Range("F115:V115").Value2
$rowB = $ws.Range("F116:V116").Value2
$ws.Range("F115:V115").Value2 = $rowB
$ws.Range("F116:V116").Value2 = $rowA

$rowA = $ws.Range("F119:V119").Value2
$rowB = $ws.Range("F120:V120").Value2
$ws.Range("F119:V119").Value2 = $rowB
$ws.Range("F120:V120").Value2 = $rowA

$rowA = $ws.Range("F123:V123").Value2
$rowB = $ws.Range("F124:V124").Value2
$ws.Range("F123:V123").Value2 = $rowB
$ws.Range("F124:V124").Value2 = $rowA

$rowA = $ws.Range("F127:V127").Value2
$rowB = $ws.Range("F128:V128").Value2
$ws.Range("F127:V127").Value2 = $rowB
$ws.Range("F128:V128").Value2 = $rowA

$rowA = $ws.Range("F147:V147").Value2
$rowB = $ws.Range("F148:V148").Value2
$ws.Range("F147:V147").Value2 = $rowB
$ws.Range("F148:V148").Value2 = $rowA

# --- Append two new match rows (155, 156) ---
# Copy row-154 formatting pattern (bold/border on A, date format on E) onto the new rows
$ws.Range("A154:V154").Copy()
$ws.Range("A155:V155").PasteSpecial(-4122)
$ws.Range("A154:V154").Copy()
$ws.Range("A156:V156").PasteSpecial(-4122)

$ws.Cells.Item(155, 1).Value2 = 154
$ws.Cells.Item(155, 2).Value2 = "turkey"
$ws.Cells.Item(155, 3).Value2 = "super-lig"
$ws.Cells.Item(155, 4).Value2 = "2023-2024"
$ws.Cells.Item(155, 5).Value2 = 45280.75
$ws.Cells.Item(155, 6).Value2 = "Antalyaspor"
$ws.Cells.Item(155, 7).Value2 = 0
$ws.Cells.Item(155, 8).Value2 = "Kasimpasa"
$ws.Cells.Item(155, 9).Value2 = 0
$ws.Cells.Item(155, 10).Value2 = 1.74
$ws.Cells.Item(155, 11).Value2 = "14/12/2023 09:42"
$ws.Cells.Item(155, 12).Value2 = 2.01
$ws.Cells.Item(155, 13).Value2 = "20/12/2023 17:57"
$ws.Cells.Item(155, 14).Value2 = 3.98
$ws.Cells.Item(155, 15).Value2 = "14/12/2023 09:42"
$ws.Cells.Item(155, 16).Value2 = 3.75
$ws.Cells.Item(155, 17).Value2 = "20/12/2023 17:57"
$ws.Cells.Item(155, 18).Value2 = 4.61
$ws.Cells.Item(155, 19).Value2 = "14/12/2023 09:42"
$ws.Cells.Item(155, 20).Value2 = 3.77
$ws.Cells.Item(155, 21).Value2 = "20/12/2023 17:57"
$ws.Cells.Item(155, 22).Value2 = "https://www.betexplorer.com/football/turkey/super-lig/antalyaspor-kasimpasa/neFLSzeQ/"

$ws.Cells.Item(156, 1).Value2 = 155
$ws.Cells.Item(156, 2).Value2 = "turkey"
$ws.Cells.Item(156, 3).Value2 = "super-lig"
$ws.Cells.Item(156, 4).Value2 = "2023-2024"
$ws.Cells.Item(156, 5).Value2 = 45280.75
$ws.Cells.Item(156, 6).Value2 = "Galatasaray"
$ws.Cells.Item(156, 7).Value2 = 1
$ws.Cells.Item(156, 8).Value2 = "Karagumruk"
$ws.Cells.Item(156, 9).Value2 = 0
$ws.Cells.Item(156, 10).Value2 = 1.23
$ws.Cells.Item(156, 11).Value2 = "14/12/2023 09:42"
$ws.Cells.Item(156, 12).Value2 = 1.27
$ws.Cells.Item(156, 13).Value2 = "20/12/2023 17:53"
$ws.Cells.Item(156, 14).Value2 = 6.94
$ws.Cells.Item(156, 15).Value2 = "14/12/2023 09:42"
$ws.Cells.Item(156, 16).Value2 = 6.61
$ws.Cells.Item(156, 17).Value2 = "20/12/2023 17:59"
$ws.Cells.Item(156, 18).Value2 = 10.54
$ws.Cells.Item(156, 19).Value2 = "14/12/2023 09:42"
$ws.Cells.Item(156, 20).Value2 = 10.32
$ws.Cells.Item(156, 21).Value2 = "20/12/2023 17:59"
$ws.Cells.Item(156, 22).Value2 = "https://www.betexplorer.com/football/turkey/super-lig/galatasaray-f-karagumruk/02ZQ6gY6/"
